$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.115.72'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.45%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.925.88'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '330.87'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4725'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4055'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.96%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '53.04'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08428'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.50%  '
$ws.Range('E11').Value = '  -5.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.29'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.93%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.927.53'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.504'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.35%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.094'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -5.93%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '90.61'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001066'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06574'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.13'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.750'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.106.32'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.40'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.287'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.132.28'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '154.09'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.09'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.156'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.17%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.713'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -10.04%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '123.68'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9777'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.15%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09609'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.67%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.451'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.555'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.636'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '9.031'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02313'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -4.96%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06179'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.33%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.234'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -7.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6165'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.07'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.002'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1904'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.34%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.294'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5881'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.66%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.86'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.93%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.030'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -7.39%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06829'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '110.09'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.01%  '
